# Update gh-pages output data (attendance numbers) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# 展览 sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 445
$wsExpo.Range("F5").Value = 482
$wsExpo.Range("F6").Value = 277
$wsExpo.Range("F7").Value = 2527
$wsExpo.Range("F9").Value = 6774

# 全部类型 sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 445
$wsAll.Range("F5").Value = 482
$wsAll.Range("F6").Value = 277
$wsAll.Range("F9").Value = 2527
$wsAll.Range("F11").Value = 6774
